# Applies the weekly update: a new price record is inserted at row 172,
# shifting the existing rows 172-193 down to 173-194 (the previous last
# row, 193, becomes row 194 and keeps its original values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 172; this shifts rows 172:193 down to 173:194
$ws.Rows.Item(172).Insert()

# Populate the new row 172 with a copy of what is now row 173 (the old row 172),
# then overwrite the four fields that actually differ for the new record
# (Fecha, Variedad, Volumen, Origen).
$ws.Range("A172").Value = 9
$ws.Range("B172").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C172").Value = "Metropolitana"
$ws.Range("D172").Value = 45218
$ws.Range("D172").NumberFormat = $ws.Range("D173").NumberFormat
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = 100112022
$ws.Range("G172").Value = "Arveja Verde"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 70
$ws.Range("K172").Value = 26000
$ws.Range("L172").Value = 28000
$ws.Range("M172").Value = 27000
$ws.Range("N172").Value = '$/malla 25 kilos'
$ws.Range("O172").Value = "Provincia de Limarí"
$ws.Range("P172").Value = 1080
$ws.Range("Q172").Value = 25
$ws.Range("R172").Value = "Hortaliza"
